$wb = $excel.ActiveWorkbook

$wsUsers      = $wb.Worksheets.Item("Users")
$wsAppName    = $wb.Worksheets.Item("AppName")
$wsModuleName = $wb.Worksheets.Item("ModuleName")
$wsGiftLog    = $wb.Worksheets.Item("GiftLog")
$wsGiftEdit   = $wb.Worksheets.Item("GiftEdit")

# --- Users sheet: the submitted-for test contact was renamed, and the cell
#     now wraps its text (new style) ---
$wsUsers.Range("A2").Value = "Julie Carthane"
$wsUsers.Range("A2").WrapText = $true

# --- GiftLog sheet: same contact rename in the "SubmittedFor" column, with
#     the same wrap-text styling, so row 2 grows tall enough to show 2 lines ---
$wsGiftLog.Range("B2").Value = "Julie Carthane"
$wsGiftLog.Range("B2").WrapText = $true
$wsGiftLog.Rows.Item(2).RowHeight = 30

# --- Remembered cursor/selection position on a couple of sheets ---
[void]$wsGiftLog.Range("D18").Select()
[void]$wsModuleName.Range("E14").Select()

# --- Window geometry (best effort; cosmetic window chrome) ---
$win = $wb.Windows.Item(1)
$win.Left = 22932
$win.Top = -108
$win.Width = 23256
$win.Height = 12456

# ModuleName becomes the active / selected tab - activate it last so it is
# the one persisted as tabSelected / workbook activeTab.
[void]$wsModuleName.Activate()
